$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "RPM" column header in L3, matching header style
$ws.Range("L3").Value = "RPM"
$ws.Range("L3").Font.Bold = $true

# RPM values for each engine row
$ws.Range("L4").Value = 9800
$ws.Range("L5").Value = 9800
$ws.Range("L6").Value = 8000
$ws.Range("L7").Value = 9000
$ws.Range("L8").Value = 8500
$ws.Range("L9").Value = 7500
$ws.Range("L10").Value = 7500
$ws.Range("L11").Value = 9000
$ws.Range("L12").Value = 9000

# Update the view: select J22, clear the previous frozen/top-left scroll position
$ws.Range("J22").Select()

# Approximate the workbook window geometry changes (best effort)
$win = $excel.ActiveWindow
$win.Left = 0
$win.Top = 0
$win.Width = 25600
$win.Height = 14460
